# Update column B (Taxonsorteringsordning) values for rows 2-17:
# each numeric value increases by 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 4
    }
}
